# Controle.xlsx -- "Processo Externo e atualização da Planilha de Controle"
#
# 1) Kanban board (sheet "Kanban"): several cards move between the
#    To Do / Doing / Done columns (the "X" marker cell shifts column),
#    and the stale note in G4 ("Leitura OK. Erro na gravação") is cleared.
# 2) Bugs sheet: a new bug row is logged.
# 3) The active sheet/selection moves from "Bugs" back to "Kanban".

$wb = $excel.ActiveWorkbook

$kanban = $wb.Worksheets.Item("Kanban")
$bugs   = $wb.Worksheets.Item("Bugs")

# --- Kanban: move cards between columns -----------------------------
# Row 4: "Solicitante - Relacionamento com Entidades de Facturamento"
#   Doing (D) -> Done (E); also clear the stale note in G4.
$kanban.Range("D4").Clear() | Out-Null
$kanban.Range("E4").Value = "X"
$kanban.Range("G4").ClearContents() | Out-Null

# Row 5: "Processo Externo - Endereços - Otimizações"
#   Doing (D) -> Done (E)
$kanban.Range("D5").Clear() | Out-Null
$kanban.Range("E5").Value = "X"

# Row 6: "Solicitante - Endereço"
#   To Do (C) -> Done (E)
$kanban.Range("C6").Clear() | Out-Null
$kanban.Range("E6").Value = "X"

# Row 8: "Processo Externo - Adequações"
#   To Do (C) -> Done (E)
$kanban.Range("C8").Clear() | Out-Null
$kanban.Range("E8").Value = "X"

# Row 9: "Processo Interno - Adequações"
#   To Do (C) -> Doing (D)
$kanban.Range("C9").Clear() | Out-Null
$kanban.Range("D9").Value = "X"

# Row 12: "Otimização da busca de Endereço por CEP"
#   To Do (C) -> Done (E)
$kanban.Range("C12").Clear() | Out-Null
$kanban.Range("E12").Value = "X"

# --- Bugs: log a new bug ---------------------------------------------
$bugs.Range("A5").Value = 4
$bugs.Range("B5").Value = "Alta"
$bugs.Range("C5").Value = "Verificar data e hora do Histórico: está exibindo 1 dia a mais e não está gravando a hora"
$bugs.Range("D5").Value = "Em aberto"

# --- Selection / active sheet -----------------------------------------
# Bugs used to be the active tab with D5 selected; now Kanban is active
# with D9 selected, and Bugs keeps a plain A6 selection.
$bugs.Range("A6").Select() | Out-Null
$kanban.Activate() | Out-Null
$kanban.Range("D9").Select() | Out-Null
